$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[R] GA(Genetic Algorithm-유전 알고리즘) 파라미터에 따른 Runtime Test"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/R-GAGenetic-Algorithm-%EC%9C%A0%EC%A0%84-%EC%95%8C%EA%B3%A0%EB%A6%AC%EC%A6%98-%ED%8C%8C%EB%9D%BC%EB%AF%B8%ED%84%B0%EC%97%90-%EB%94%B0%EB%A5%B8-Runtime-Test-1"

$ws.Range("D9").Value = "[공지] SIAI 설립 관련 서류"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/siai-foundation-letter/#utm_source=rss&utm_medium=rss&utm_campaign=siai-foundation-letter"

$ws.Range("D45").Value = "MCMC (Markov Chain Monte Carlo)"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/97"

$ws.Range("D52").Value = "서울특별시 시간별 (초)미세먼지 03: 극단값 확인"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2511&utm_source=rss&utm_medium=rss&utm_campaign=%25ec%2584%259c%25ec%259a%25b8%25ed%258a%25b9%25eb%25b3%2584%25ec%258b%259c-%25ec%258b%259c%25ea%25b0%2584%25eb%25b3%2584-%25ec%25b4%2588%25eb%25af%25b8%25ec%2584%25b8%25eb%25a8%25bc%25ec%25a7%2580-03-%25ea%25b7%25b9%25eb%258b%25a8%25ea%25b0%2592-%25ed%2599%2595%25ec%259d%25b8"
